# "01/02 - Commit - Final D03 Changes :)"
# Append four new Variable/Value rows (14-17) to the "Environments_OnGoing"
# sheet describing the new DirectSales Products/Files related-list links,
# and make that sheet the active/selected one (it was previously
# "D03NonQuotableProducts" that was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Environments_OnGoing")

# Row 14
$ws.Range("A14").Value = "DirectSalesProductsOpportunity"
$ws.Range("B14").Value = "/lightning/r/OpportunityLineItem/"

# Row 15 (value entered before the variable name, matching original authoring order)
$ws.Range("B15").Value = "/related/OpportunityLineItems/view"
$ws.Range("A15").Value = "DirectSalesProductsRelatedView"

# Row 16
$ws.Range("A16").Value = "DirectSalesFiles"
$ws.Range("B16").Value = "/lightning/r/AttachedContentDocument/"

# Row 17
$ws.Range("A17").Value = "DirectSalesFilesRelatedView"
$ws.Range("B17").Value = "/related/AttachedContentDocuments/view"

# Make this sheet the active tab (moving away from D03NonQuotableProducts),
# and leave selection on the next empty row as the author did.
$ws.Activate()
$ws.Range("B18").Select()
